$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 and Row 3 fully swap their species-record content ---
# Save current (pre-edit) values of row 2 and row 3 so we can cross-assign them.
$A2 = $ws.Range("A2").Value()
$B2 = $ws.Range("B2").Value()
$C2 = $ws.Range("C2").Value()
$D2 = $ws.Range("D2").Value()
$E2 = $ws.Range("E2").Value()
$F2 = $ws.Range("F2").Value()
$G2 = $ws.Range("G2").Value()
$H2 = $ws.Range("H2").Value()

$A3 = $ws.Range("A3").Value()
$B3 = $ws.Range("B3").Value()
$C3 = $ws.Range("C3").Value()
$D3 = $ws.Range("D3").Value()
$E3 = $ws.Range("E3").Value()
$F3 = $ws.Range("F3").Value()
$G3 = $ws.Range("G3").Value()
$H3 = $ws.Range("H3").Value()

# Row 2 becomes the "Tretåig hackspett" record (previously on row 3),
# and gains the Aktivitet note "äldre spår" (column M) plus the
# placeholder Ålder-Stadium / Kön / Metod cells (K, L, N).
$ws.Range("A2").Value = 111396308
$ws.Range("B2").Value = $B3
$ws.Range("C2").Value = $C3
$ws.Range("D2").Value = $D3
$ws.Range("E2").Value = $E3
$ws.Range("F2").Value = $F3
$ws.Range("G2").Value = $G3
$ws.Range("H2").Value = $H3
$ws.Range("K2").Font.Bold = $false
$ws.Range("L2").Font.Bold = $false
$ws.Range("M2").Value = "äldre spår"
$ws.Range("N2").Font.Bold = $false
$ws.Range("Q2").Value = 625151.1577179903
$ws.Range("R2").Value = 7209567.512248591

# Row 3 becomes a plain "Knärot" record (matching the other rows),
# and loses the Ålder-Stadium / Kön / Aktivitet / Metod cells it had.
$ws.Range("A3").Value = 111396322
$ws.Range("B3").Value = $B2
$ws.Range("C3").Value = $C2
$ws.Range("D3").Value = $D2
$ws.Range("E3").Value = $E2
$ws.Range("F3").Value = $F2
$ws.Range("G3").Value = $G2
$ws.Range("H3").Value = $H2
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("Q3").Value = 625269.4478252844
$ws.Range("R3").Value = 7209630.115016816

# --- Rows 4-19: the Id (A) and coordinate (Q, R) values were reshuffled
#     among the rows; everything else about each record is unchanged. ---
$ws.Range("A4").Value = 111396326
$ws.Range("Q4").Value = 625397.1584455093
$ws.Range("R4").Value = 7209589.718691397

$ws.Range("A5").Value = 111396321
$ws.Range("Q5").Value = 625240.2002264742
$ws.Range("R5").Value = 7209649.650274927

$ws.Range("A6").Value = 111396319
$ws.Range("Q6").Value = 625228.8129008666
$ws.Range("R6").Value = 7209607.642547456

$ws.Range("A7").Value = 111396312
$ws.Range("Q7").Value = 625242.7087276473
$ws.Range("R7").Value = 7209468.80281719

$ws.Range("A8").Value = 111396315
$ws.Range("Q8").Value = 625167.9685939638
$ws.Range("R8").Value = 7209530.9258211

$ws.Range("A9").Value = 111396311
$ws.Range("Q9").Value = 625271.0561409625
$ws.Range("R9").Value = 7209511.101565193

$ws.Range("A10").Value = 111396325
$ws.Range("Q10").Value = 625389.9085714296
$ws.Range("R10").Value = 7209580.514361567

$ws.Range("A11").Value = 111396310
$ws.Range("Q11").Value = 625289.0018867656
$ws.Range("R11").Value = 7209518.212698339

$ws.Range("A12").Value = 111396317
$ws.Range("Q12").Value = 625153.5624699651
$ws.Range("R12").Value = 7209550.662191558

$ws.Range("A13").Value = 111396323
$ws.Range("Q13").Value = 625301.6605433678
$ws.Range("R13").Value = 7209610.70454926

$ws.Range("A14").Value = 111396318
$ws.Range("Q14").Value = 625177.6865340136
$ws.Range("R14").Value = 7209552.099144561

$ws.Range("A15").Value = 111396314
$ws.Range("Q15").Value = 625202.8383709632
$ws.Range("R15").Value = 7209539.171001118

$ws.Range("A16").Value = 111396309
$ws.Range("Q16").Value = 625341.71034419
$ws.Range("R16").Value = 7209536.108963673

$ws.Range("A17").Value = 111396316
$ws.Range("Q17").Value = 625153.7279882778
$ws.Range("R17").Value = 7209526.513740451

$ws.Range("A18").Value = 111396324
$ws.Range("Q18").Value = 625335.6676841485
$ws.Range("R18").Value = 7209609.168182318

$ws.Range("A19").Value = 111396313
$ws.Range("Q19").Value = 625231.5510770321
$ws.Range("R19").Value = 7209481.895207534
